$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 16668500
$ws.Range("J70").Value = 1999.7778
$ws.Range("L70").Value = 5999.3334
$ws.Range("N70").Value = -6539.3334
$ws.Range("H73").Value = 16668500
$ws.Range("J73").Value = 1999.7778
$ws.Range("L73").Value = 5999.3334
$ws.Range("N73").Value = -7871.3334
$ws.Range("H98").Value = 27633.154
$ws.Range("J98").Value = 15000
$ws.Range("L98").Value = 15000
$ws.Range("N98").Value = -17996
$ws.Range("H100").Value = 22864298
$ws.Range("I100").Value = 35859040
$ws.Range("K100").Value = 35859040
$ws.Range("M100").Value = -35858499
$ws.Range("H112").Value = 2726.8823
$ws.Range("I112").Value = 3511
$ws.Range("J112").Value = 2485.6155
$ws.Range("K112").Value = 10533
$ws.Range("L112").Value = 7456.8465
$ws.Range("M112").Value = -9425
$ws.Range("N112").Value = -9672.8465
$ws.Range("H122").Value = 27633.154
$ws.Range("J122").Value = 15000
$ws.Range("L122").Value = 45000
$ws.Range("N122").Value = -49900
$ws.Range("H125").Value = 8427.286
$ws.Range("I125").Value = 17998
$ws.Range("J125").Value = 4599
$ws.Range("K125").Value = 161982
$ws.Range("L125").Value = 41391
$ws.Range("M125").Value = -159522
$ws.Range("N125").Value = -46311
$ws.Range("H132").Value = 1963900.6
$ws.Range("I132").Value = 3027.422
$ws.Range("K132").Value = 9082.266
$ws.Range("M132").Value = -6552.266
$ws.Range("H138").Value = 293987.34
$ws.Range("I138").Value = 439361.1
$ws.Range("J138").Value = 3239.8
$ws.Range("K138").Value = 1318083.3
$ws.Range("L138").Value = 9719.400000000001
$ws.Range("M138").Value = -1312943.3
$ws.Range("N138").Value = -19999.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6598.3174
$ws.Range("I32").Value = 6503.509
$ws.Range("J32").Value = 7499
$ws.Range("K32").Value = 6503.509
$ws.Range("L32").Value = 7499
$ws.Range("M32").Value = -6216.509
$ws.Range("N32").Value = -8073
$ws.Range("H61").Value = 8727.966
$ws.Range("I61").Value = 10256.5
$ws.Range("K61").Value = 10256.5
$ws.Range("M61").Value = -10044.5
$ws.Range("H74").Value = 13156.96
$ws.Range("I74").Value = 18551.5
$ws.Range("J74").Value = 3566.6667
$ws.Range("K74").Value = 18551.5
$ws.Range("L74").Value = 3566.6667
$ws.Range("M74").Value = -17677.5
$ws.Range("N74").Value = -5314.6667
$ws.Range("H77").Value = 13156.96
$ws.Range("I77").Value = 18551.5
$ws.Range("J77").Value = 3566.6667
$ws.Range("K77").Value = 92757.5
$ws.Range("L77").Value = 17833.3335
$ws.Range("M77").Value = -88389.5
$ws.Range("N77").Value = -26569.3335
$ws.Range("H102").Value = 11764.381
$ws.Range("I102").Value = 16541.143
$ws.Range("K102").Value = 16541.143
$ws.Range("M102").Value = -14919.143
$ws.Range("H110").Value = 2613.25
$ws.Range("I110").Value = 1942.6666
$ws.Range("K110").Value = 1942.6666
$ws.Range("M110").Value = 102.3334
$ws.Range("H113").Value = 112297
$ws.Range("J113").Value = 112297
$ws.Range("L113").Value = 112297
$ws.Range("N113").Value = -120975
$ws.Range("H122").Value = 941825.1
$ws.Range("I122").Value = 3762
$ws.Range("K122").Value = 11286
$ws.Range("M122").Value = -8836
$ws.Range("H132").Value = 2977.9778
$ws.Range("I132").Value = 2750.9678
$ws.Range("J132").Value = 3480.6428
$ws.Range("K132").Value = 8252.903399999999
$ws.Range("L132").Value = 10441.9284
$ws.Range("M132").Value = -5722.903399999999
$ws.Range("N132").Value = -15501.9284
$ws.Range("H136").Value = 8727.966
$ws.Range("I136").Value = 10256.5
$ws.Range("K136").Value = 30769.5
$ws.Range("M136").Value = -28219.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 342034850
$ws.Range("J108").Value = 342034850
$ws.Range("L108").Value = 342034850
$ws.Range("N108").Value = -342042530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1234.8
$ws.Range("I5").Value = 213
$ws.Range("J5").Value = 1916
$ws.Range("K5").Value = 213
$ws.Range("L5").Value = 1916
$ws.Range("M5").Value = -101
$ws.Range("N5").Value = -2140
$ws.Range("H134").Value = 5330.353
$ws.Range("I134").Value = 7032.864
$ws.Range("K134").Value = 21098.592
$ws.Range("M134").Value = -18563.592

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 998.5
$ws.Range("I33").Value = 998.5
$ws.Range("K33").Value = 5991
$ws.Range("M33").Value = -5708
$ws.Range("H97").Value = 55772.637
$ws.Range("I97").Value = 76575
$ws.Range("J97").Value = 299.66666
$ws.Range("K97").Value = 229725
$ws.Range("L97").Value = 898.9999799999999
$ws.Range("M97").Value = -229229
$ws.Range("N97").Value = -1890.99998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11071.652
$ws.Range("I122").Value = 8179.8335
$ws.Range("J122").Value = 14226.363
$ws.Range("K122").Value = 24539.5005
$ws.Range("L122").Value = 42679.089
$ws.Range("M122").Value = -22089.5005
$ws.Range("N122").Value = -47579.089
$ws.Range("H132").Value = 4360.25
$ws.Range("I132").Value = 5143.3213
$ws.Range("K132").Value = 15429.9639
$ws.Range("M132").Value = -12899.9639

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 23871.54
$ws.Range("J40").Value = 15665
$ws.Range("L40").Value = 15665
$ws.Range("N40").Value = -15937
$ws.Range("H46").Value = 2113658.5
$ws.Range("I46").Value = 768.5333000000001
$ws.Range("J46").Value = 3874400
$ws.Range("K46").Value = 768.5333000000001
$ws.Range("L46").Value = 3874400
$ws.Range("M46").Value = -580.5333000000001
$ws.Range("N46").Value = -3874776
$ws.Range("H61").Value = 3809.55
$ws.Range("I61").Value = 806.06665
$ws.Range("K61").Value = 806.06665
$ws.Range("M61").Value = -604.06665
$ws.Range("H68").Value = 3082.1667
$ws.Range("I68").Value = 2538.6
$ws.Range("J68").Value = 5800
$ws.Range("K68").Value = 2538.6
$ws.Range("L68").Value = 5800
$ws.Range("M68").Value = -1789.6
$ws.Range("N68").Value = -7298
$ws.Range("H71").Value = 3082.1667
$ws.Range("I71").Value = 2538.6
$ws.Range("J71").Value = 5800
$ws.Range("K71").Value = 12693
$ws.Range("L71").Value = 29000
$ws.Range("M71").Value = -8949
$ws.Range("N71").Value = -36488
$ws.Range("H82").Value = 2947.8462
$ws.Range("I82").Value = 3103.5557
$ws.Range("J82").Value = 2597.5
$ws.Range("K82").Value = 3103.5557
$ws.Range("L82").Value = 2597.5
$ws.Range("M82").Value = -2742.5557
$ws.Range("N82").Value = -3319.5
$ws.Range("H85").Value = 2947.8462
$ws.Range("I85").Value = 3103.5557
$ws.Range("J85").Value = 2597.5
$ws.Range("K85").Value = 3103.5557
$ws.Range("L85").Value = 2597.5
$ws.Range("M85").Value = -1855.5557
$ws.Range("N85").Value = -5093.5
$ws.Range("H113").Value = 3809.55
$ws.Range("I113").Value = 806.06665
$ws.Range("K113").Value = 806.06665
$ws.Range("M113").Value = 1363.93335
$ws.Range("H122").Value = 5191.0347
$ws.Range("I122").Value = 5157.6
$ws.Range("K122").Value = 15472.8
$ws.Range("M122").Value = -13022.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 491154.84
$ws.Range("I62").Value = 571439
$ws.Range("J62").Value = 9450
$ws.Range("K62").Value = 571439
$ws.Range("L62").Value = 9450
$ws.Range("M62").Value = -570815
$ws.Range("N62").Value = -10698
$ws.Range("H65").Value = 491154.84
$ws.Range("I65").Value = 571439
$ws.Range("J65").Value = 9450
$ws.Range("K65").Value = 2857195
$ws.Range("L65").Value = 47250
$ws.Range("M65").Value = -2854075
$ws.Range("N65").Value = -53490
$ws.Range("H107").Value = 24824.385
$ws.Range("I107").Value = 1893.0834
$ws.Range("K107").Value = 5679.2502
$ws.Range("M107").Value = -3759.2502
$ws.Range("H132").Value = 17314.074
$ws.Range("I132").Value = 27889.133
$ws.Range("K132").Value = 83667.399
$ws.Range("M132").Value = -81137.399
$ws.Range("H136").Value = 683214.2
$ws.Range("I136").Value = 1038555.06
$ws.Range("J136").Value = 16950
$ws.Range("K136").Value = 3115665.18
$ws.Range("M136").Value = -3113115.18
$ws.Range("N136").Value = -55950
